$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused | Clear Glass Lens
$ws.Range("H33").Value = 225.96297
$ws.Range("I33").Value = 245.78261
$ws.Range("K33").Value = 245.78261
$ws.Range("M33").Value = -16.78261000000001

# Row 39: Riches' Brew | Hi-Potion of Mind
$ws.Range("H39").Value = 439.46155
$ws.Range("I39").Value = 262.57144
$ws.Range("J39").Value = 645.8333
$ws.Range("K39").Value = 787.71432
$ws.Range("L39").Value = 1937.4999
$ws.Range("M39").Value = -491.71432
$ws.Range("N39").Value = -2529.4999

# Row 55: A Real Smooth Move | Lanolin
$ws.Range("H55").Value = 471.42856
$ws.Range("I55").Value = 301
$ws.Range("J55").Value = 499.83334
$ws.Range("K55").Value = 301
$ws.Range("L55").Value = 499.83334
$ws.Range("M55").Value = -87
$ws.Range("N55").Value = -927.83334

# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 1874.1
$ws.Range("I112").Value = 717.6
$ws.Range("J112").Value = 2259.6
$ws.Range("K112").Value = 2152.8
$ws.Range("L112").Value = 6778.799999999999
$ws.Range("M112").Value = -1044.8
$ws.Range("N112").Value = -8994.799999999999

# Row 125: Body over Mind | Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 71746.8
$ws.Range("I125").Value = 170271.67
$ws.Range("K125").Value = 1532445.03
$ws.Range("M125").Value = -1529985.03

# Row 129: Practical Command | Commanding Craftsman's Draught
$ws.Range("H129").Value = 1260.3658
$ws.Range("J129").Value = 1291.9231
$ws.Range("L129").Value = 3875.7693
$ws.Range("N129").Value = -13875.7693

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 16710.613
$ws.Range("J32").Value = 4185.5
$ws.Range("L32").Value = 4185.5
$ws.Range("N32").Value = -4759.5

# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 1700.75
$ws.Range("I45").Value = 1499
$ws.Range("J45").Value = 1768
$ws.Range("K45").Value = 1499
$ws.Range("L45").Value = 1768
$ws.Range("M45").Value = -1122
$ws.Range("N45").Value = -2522

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 1520.5862
$ws.Range("I61").Value = 1430.66
$ws.Range("K61").Value = 1430.66
$ws.Range("M61").Value = -1218.66

# Row 97: Ore for Me | High Steel Ingot
$ws.Range("H97").Value = 824.25
$ws.Range("I97").Value = 798.5
$ws.Range("K97").Value = 798.5
$ws.Range("M97").Value = -302.5

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 4362.8335
$ws.Range("I122").Value = 4781.75
$ws.Range("J122").Value = 3525
$ws.Range("K122").Value = 14345.25
$ws.Range("L122").Value = 10575
$ws.Range("M122").Value = -11895.25
$ws.Range("N122").Value = -15475

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1520.5862
$ws.Range("I136").Value = 1430.66
$ws.Range("K136").Value = 4291.98
$ws.Range("M136").Value = -1741.98

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 4764828.5
$ws.Range("I105").Value = 8931652
$ws.Range("J105").Value = 2743.9285
$ws.Range("K105").Value = 8931652
$ws.Range("L105").Value = 2743.9285
$ws.Range("M105").Value = -8929905
$ws.Range("N105").Value = -6237.9285

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 30279.611
$ws.Range("I107").Value = 40978.617
$ws.Range("J107").Value = 2462.2
$ws.Range("K107").Value = 40978.617
$ws.Range("L107").Value = 2462.2
$ws.Range("M107").Value = -39058.617
$ws.Range("N107").Value = -6302.2

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall | Elm Lumber
$ws.Range("H22").Value = 407.5
$ws.Range("I22").Value = 244.84616
$ws.Range("J22").Value = 830.4
$ws.Range("K22").Value = 244.84616
$ws.Range("L22").Value = 830.4
$ws.Range("M22").Value = 105.15384
$ws.Range("N22").Value = -1530.4

# Row 62: Splinter in the Sewers | Cedar Lumber
$ws.Range("H62").Value = 86084.164
$ws.Range("I62").Value = 127126.25
$ws.Range("K62").Value = 127126.25
$ws.Range("M62").Value = -126502.25

# Row 65: The Lumber of Their Discontent (L) | Cedar Lumber
$ws.Range("H65").Value = 86084.164
$ws.Range("I65").Value = 127126.25
$ws.Range("K65").Value = 635631.25
$ws.Range("M65").Value = -632511.25

# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 2343.889
$ws.Range("I99").Value = 2283.3333
$ws.Range("J99").Value = 2465
$ws.Range("K99").Value = 2283.3333
$ws.Range("L99").Value = 2465
$ws.Range("M99").Value = -785.3332999999998
$ws.Range("N99").Value = -5461

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 2343.889
$ws.Range("I126").Value = 2283.3333
$ws.Range("J126").Value = 2465
$ws.Range("K126").Value = 6849.999899999999
$ws.Range("L126").Value = 7395
$ws.Range("M126").Value = -4379.999899999999
$ws.Range("N126").Value = -12335

$ws = $wb.Worksheets.Item("CUL")
# Row 140: Sweet, Sweet Bean Juice | Mesquite Juice
$ws.Range("H140").Value = 1827.6552
$ws.Range("I140").Value = 939.5
$ws.Range("J140").Value = 4619
$ws.Range("K140").Value = 2818.5
$ws.Range("L140").Value = 13857
$ws.Range("M140").Value = 2361.5
$ws.Range("N140").Value = -24217

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 4155.4443
$ws.Range("I102").Value = 4100
$ws.Range("J102").Value = 4224.75
$ws.Range("K102").Value = 4100
$ws.Range("L102").Value = 4224.75
$ws.Range("M102").Value = -2478
$ws.Range("N102").Value = -7468.75

# Row 109: You're My Wonderhall | Hematite Earrings of Healing
$ws.Range("H109").Value = 9272
$ws.Range("J109").Value = 9272
$ws.Range("L109").Value = 9272
$ws.Range("N109").Value = -11352

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 3667.8333
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 3728.5454
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 11185.6362
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -16085.6362

# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 2975.4211
$ws.Range("I126").Value = 2771.7693
$ws.Range("J126").Value = 3416.6667
$ws.Range("K126").Value = 8315.3079
$ws.Range("L126").Value = 10250.0001
$ws.Range("M126").Value = -5845.3079
$ws.Range("N126").Value = -15190.0001

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 2591.4
$ws.Range("I132").Value = 2061.182
$ws.Range("J132").Value = 4049.5
$ws.Range("K132").Value = 6183.545999999999
$ws.Range("L132").Value = 12148.5
$ws.Range("M132").Value = -3653.545999999999
$ws.Range("N132").Value = -17208.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 3372.9473
$ws.Range("I7").Value = 2962.6
$ws.Range("J7").Value = 3828.889
$ws.Range("K7").Value = 2962.6
$ws.Range("L7").Value = 3828.889
$ws.Range("M7").Value = -2850.6
$ws.Range("N7").Value = -4052.889

# Row 68: You Could Say It's a Moving Target | Wyvern Leather
$ws.Range("H68").Value = 3185.5715
$ws.Range("I68").Value = 2669.2
$ws.Range("J68").Value = 4476.5
$ws.Range("K68").Value = 2669.2
$ws.Range("L68").Value = 4476.5
$ws.Range("M68").Value = -1920.2
$ws.Range("N68").Value = -5974.5

# Row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Range("H71").Value = 3185.5715
$ws.Range("I71").Value = 2669.2
$ws.Range("J71").Value = 4476.5
$ws.Range("K71").Value = 13346
$ws.Range("L71").Value = 22382.5
$ws.Range("M71").Value = -9602
$ws.Range("N71").Value = -29870.5

# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 100004000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 3372.9473
$ws.Range("I126").Value = 2962.6
$ws.Range("J126").Value = 3828.889
$ws.Range("K126").Value = 8887.799999999999
$ws.Range("L126").Value = 11486.667
$ws.Range("M126").Value = -6417.799999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display | Ruby Cotton Cloth
$ws.Range("H96").Value = 909.75
$ws.Range("I96").Value = 651
$ws.Range("J96").Value = 1065
$ws.Range("K96").Value = 651
$ws.Range("L96").Value = 1065
$ws.Range("M96").Value = 722
$ws.Range("N96").Value = -3811

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 5518.8
$ws.Range("I126").Value = 6411.125
$ws.Range("K126").Value = 19233.375
$ws.Range("M126").Value = -16763.375
